# Rearranges species-observation rows 26, 28, 29, 30, 31 on the active sheet.
# Row 27 is untouched. The affected columns are A, B, D, E, F, G, H, Q, R.
#
# The net effect (per the target diff) is:
#   new row 26 <= old row 28
#   new row 28 <= old row 29
#   new row 29 <= old row 26
#   new row 30 <= old row 31
#   new row 31 <= old row 30

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

function Get-RowValues($ws, $row, $cols) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Range("$c$row").Value2
    }
    return $vals
}

function Set-RowValues($ws, $row, $cols, $vals) {
    foreach ($c in $cols) {
        $ws.Range("$c$row").Value2 = $vals[$c]
    }
}

# Snapshot the current ("before") contents of every row involved.
$row26 = Get-RowValues $ws 26 $cols
$row28 = Get-RowValues $ws 28 $cols
$row29 = Get-RowValues $ws 29 $cols
$row30 = Get-RowValues $ws 30 $cols
$row31 = Get-RowValues $ws 31 $cols

# Write back the permuted values.
Set-RowValues $ws 26 $cols $row28
Set-RowValues $ws 28 $cols $row29
Set-RowValues $ws 29 $cols $row26
Set-RowValues $ws 30 $cols $row31
Set-RowValues $ws 31 $cols $row30
